$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 108-110: row-height-only tweaks (content/style unchanged) ---
$ws.Rows.Item(108).RowHeight = 41.25
$ws.Rows.Item(109).RowHeight = 41.25
$ws.Rows.Item(110).RowHeight = 27.75

# --- Row 111 ---
$ws.Range("A111").Value = 'S13'
$ws.Range("B111").Value = 'G02'
$ws.Range("C111").Value = 'Branding and logo integration'
$ws.Range("D111").Value = 'S13_G02_TF001'
$ws.Range("E111").Value = 'Add SigmaTrader logo to the app shell (AppBar/sidebar) using assets from frontend/public and ensure it fits both light and dark themes.'
$ws.Range("G111").Value = 'implemented'
$ws.Range("H111").Value = 'Initial focus on placing a small logo mark next to the SigmaTrader title in the top bar or sidebar without changing layout structure.'
$ws.Range("I111").Value = 'Decide final logo placement and sizing, then optionally extend branding to the auth/landing page hero area.'
$ws.Rows.Item(111).WrapText = $true
$ws.Rows.Item(111).VerticalAlignment = -4160
$ws.Rows.Item(111).RowHeight = 41.75

# --- Row 112 ---
$ws.Range("A112").Value = 'S14'
$ws.Range("B112").Value = 'G01'
$ws.Range("C112").Value = 'Advanced order types and stop-loss controls'
$ws.Range("D112").Value = 'S14_G01_TB001'
$ws.Range("E112").Value = 'Extend Order model and APIs to support Zerodha order types MARKET / LIMIT / SL / SL-M plus trigger price and trigger percent fields.'
$ws.Range("F112").Value = 'trigger_percent will be interpreted relative to the Zerodha last traded price (LTP), not the current limit price, and saved alongside trigger_price.'
$ws.Range("G112").Value = 'implemented'
$ws.Range("H112").Value = 'Order model, schemas, and Zerodha client now support trigger_price/trigger_percent and extended order_type values (MARKET/LIMIT/SL/SL-M).'
$ws.Range("I112").Value = 'UI still only exposes MARKET/LIMIT; S14/G03 will surface SL/SL-M and triggers in the edit dialog.'
$ws.Rows.Item(112).WrapText = $true
$ws.Rows.Item(112).VerticalAlignment = -4160
$ws.Rows.Item(112).RowHeight = 41.75

# --- Row 113 ---
$ws.Range("A113").Value = 'S14'
$ws.Range("B113").Value = 'G01'
$ws.Range("C113").Value = 'Advanced order types and stop-loss controls'
$ws.Range("D113").Value = 'S14_G01_TB002'
$ws.Range("E113").Value = 'Update execute_order to route SL and SL-M correctly to KiteConnect (trigger_price mandatory, price optional for SL-M) and add guardrails for valid stop-loss placement.'
$ws.Range("F113").Value = 'Guardrails include checking trigger_price vs LTP and direction (BUY stops below market, SELL stops above) and rejecting obviously invalid combinations with clear error messages.'
$ws.Range("G113").Value = 'implemented'
$ws.Range("H113").Value = 'execute_order applies LTP-based stop-loss guardrails and forwards trigger_price to Zerodha with appropriate pricing rules for SL vs SL-M.'
$ws.Range("I113").Value = 'Further tuning of guardrail rules can be done once real-world usage feedback is collected.'
$ws.Rows.Item(113).WrapText = $true
$ws.Rows.Item(113).VerticalAlignment = -4160
$ws.Rows.Item(113).RowHeight = 41.75

# --- Row 114 ---
$ws.Range("A114").Value = 'S14'
$ws.Range("B114").Value = 'G02'
$ws.Range("C114").Value = 'Funds and margin preview for edited orders'
$ws.Range("D114").Value = 'S14_G02_TB001'
$ws.Range("E114").Value = 'Wrap Zerodha margins and order_margins APIs in the backend and expose endpoints to fetch available funds and a margin/charges preview for a hypothetical order.'
$ws.Range("F114").Value = 'Focus first on the equity segment for Zerodha; later sprints can extend to derivatives or other brokers.'
$ws.Range("G114").Value = 'implemented'
$ws.Range("H114").Value = 'Zerodha margins and order_margins are wrapped by new /api/zerodha/margins and /api/zerodha/order-preview endpoints for funds and charges preview.'
$ws.Range("I114").Value = 'Current implementation focuses on equity segment; derivatives and other brokers can be added later.'
$ws.Rows.Item(114).WrapText = $true
$ws.Rows.Item(114).VerticalAlignment = -4160
$ws.Rows.Item(114).RowHeight = 41.75

# --- Row 115 ---
$ws.Range("A115").Value = 'S14'
$ws.Range("B115").Value = 'G02'
$ws.Range("C115").Value = 'Funds and margin preview for edited orders'
$ws.Range("D115").Value = 'S14_G02_TF002'
$ws.Range("E115").Value = 'Enhance the Waiting Queue edit dialog to display Required funds (incl. charges) and Available funds, updating dynamically as qty/price/type change.'
$ws.Range("F115").Value = 'Uses the new Zerodha preview APIs to recompute required margin after each significant edit, with warnings when funds are insufficient.'
$ws.Range("G115").Value = 'implemented'
$ws.Range("H115").Value = 'Queue edit dialog now has a Funds & charges panel that calls the new preview APIs and shows Required vs Available amounts.'
$ws.Range("I115").Value = 'Funds preview is on-demand via Recalculate; future work may auto-refresh as fields change.'
$ws.Rows.Item(115).WrapText = $true
$ws.Rows.Item(115).VerticalAlignment = -4160
$ws.Rows.Item(115).RowHeight = 41.75

# --- Row 116 ---
$ws.Range("A116").Value = 'S14'
$ws.Range("B116").Value = 'G03'
$ws.Range("C116").Value = 'Queue edit UX polish and stop-loss helpers'
$ws.Range("D116").Value = 'S14_G03_TF001'
$ws.Range("E116").Value = 'Refactor the edit queue order dialog into clear sections (Quantity & price, Stop-loss, Product & preferences, Funds) and add fields for trigger price and trigger percent.'
$ws.Range("F116").Value = 'Trigger percent will be entered as a positive or negative percentage relative to LTP and used to derive trigger_price with inline preview.'
$ws.Range("G116").Value = 'implemented'
$ws.Range("H116").Value = 'Queue edit dialog now has structured sections, BUY/SELL toggles, full order type support (MARKET/LIMIT/SL/SL-M), and complementary trigger price/percent fields with LTP-aware helpers.'
$ws.Range("I116").Value = 'Future refinements may include more granular guidance for SL placement and inline warnings when trigger levels are too close to LTP.'
$ws.Rows.Item(116).WrapText = $true
$ws.Rows.Item(116).VerticalAlignment = -4160
$ws.Rows.Item(116).RowHeight = 55.2

# --- Row 117 ---
$ws.Range("A117").Value = 'S15'
$ws.Range("B117").Value = 'G01'
$ws.Range("C117").Value = 'Zerodha GTT order support'
$ws.Range("D117").Value = 'S15_G01_TB001'
$ws.Range("E117").Value = 'Design how SigmaTrader will map queue orders and preferences into Zerodha GTT single-leg orders (trigger values, last_price source, and order payload).'
$ws.Range("F117").Value = 'GTT design should clarify when to use GTT vs regular orders, how to represent GTT status in SigmaTrader, and how TradingView alerts can request GTT creation.'
$ws.Range("G117").Value = 'pending'
$ws.Range("I117").Value = 'This design underpins safe GTT usage for CNC swing trades and must consider off-market placement and modification flows.'
$ws.Rows.Item(117).WrapText = $true
$ws.Rows.Item(117).VerticalAlignment = -4160
$ws.Rows.Item(117).RowHeight = 41.75

# --- Row 118 ---
$ws.Range("A118").Value = 'S15'
$ws.Range("B118").Value = 'G01'
$ws.Range("C118").Value = 'Zerodha GTT order support'
$ws.Range("D118").Value = 'S15_G01_TB002'
$ws.Range("E118").Value = 'Implement backend support for placing, listing, and cancelling Zerodha GTTs using KiteConnect place_gtt / get_gtts / delete_gtt, wired to per-user broker connections.'
$ws.Range("F118").Value = 'Initial focus on single-leg GTTs for equity; OCO/advanced patterns can be added later.'
$ws.Range("G118").Value = 'pending'
$ws.Range("I118").Value = 'Requires careful error handling and alignment between SigmaTrader order records and Zerodha GTT IDs.'
$ws.Rows.Item(118).WrapText = $true
$ws.Rows.Item(118).VerticalAlignment = -4160
$ws.Rows.Item(118).RowHeight = 55.2

# --- Row 119 ---
$ws.Range("A119").Value = 'S15'
$ws.Range("B119").Value = 'G01'
$ws.Range("C119").Value = 'Zerodha GTT order support'
$ws.Range("D119").Value = 'S15_G01_TF003'
$ws.Range("E119").Value = 'Extend the queue edit and manual order flows to allow creating GTT orders (instead of or in addition to regular orders) when the user selects a GTT option.'
$ws.Range("F119").Value = 'The existing "Convert to GTT" checkbox will be repurposed into a concrete GTT mode that creates or updates real GTTs at Zerodha rather than acting as a passive preference.'
$ws.Range("G119").Value = 'pending'
$ws.Range("I119").Value = 'UI should clearly distinguish between regular orders and GTTs and indicate when an order has an associated active GTT at the broker.'
$ws.Rows.Item(119).WrapText = $true
$ws.Rows.Item(119).VerticalAlignment = -4160
$ws.Rows.Item(119).RowHeight = 41.75

# --- Row 120 ---
$ws.Range("A120").Value = 'S16'
$ws.Range("B120").Value = 'G01'
$ws.Range("C120").Value = 'Paper trading architecture and configuration'
$ws.Range("D120").Value = 'S16_G01_TB001'
$ws.Range("E120").Value = 'Design paper trading execution targets (e.g., LIVE vs PAPER) and per-user/per-strategy configuration, including a selectable poll interval between 15 seconds and 4 hours.'
$ws.Range("F120").Value = 'Paper mode will reuse existing Order/Risk/Analytics models; only execution routing and price monitoring change.'
$ws.Range("G120").Value = 'pending'
$ws.Range("I120").Value = 'Document configuration in PRD and decide where poll interval is stored (global vs per-strategy vs per-user).'

# --- Row 121 ---
$ws.Range("A121").Value = 'S16'
$ws.Range("B121").Value = 'G02'
$ws.Range("C121").Value = 'Backend paper execution engine and price polling'
$ws.Range("D121").Value = 'S16_G02_TB001'
$ws.Range("E121").Value = 'Implement a paper broker engine that manages simulated orders and periodically polls Zerodha LTP for symbols with open paper orders, filling LIMIT/SL/SL-M orders when prices cross configured levels.'
$ws.Range("F121").Value = 'Price monitoring will be based on LTP polling at the user-selected interval rather than a full order book simulation.'
$ws.Range("G121").Value = 'pending'
$ws.Range("I121").Value = 'Start with simple fill rules (e.g., limit BUY fills when LTP <= limit) and extend later if needed.'

# --- Row 122 ---
$ws.Range("A122").Value = 'S16'
$ws.Range("B122").Value = 'G02'
$ws.Range("C122").Value = 'Backend paper execution engine and price polling'
$ws.Range("D122").Value = 'S16_G02_TB002'
$ws.Range("E122").Value = 'Route TradingView AUTO orders and manual queue execution through the paper engine when the strategy/user is configured for PAPER, avoiding any real Zerodha calls and marking orders as simulated.'
$ws.Range("F122").Value = 'Simulated orders will use Order.simulated = True and separate status transitions while sharing the same analytics pipeline.'
$ws.Range("G122").Value = 'pending'
$ws.Range("I122").Value = 'Ensure logging, risk checks, and analytics clearly distinguish between live and paper trades.'

# --- Row 123 ---
$ws.Range("A123").Value = 'S16'
$ws.Range("B123").Value = 'G03'
$ws.Range("C123").Value = 'Paper mode UI and analytics integration'
$ws.Range("D123").Value = 'S16_G03_TF001'
$ws.Range("E123").Value = 'Extend Settings/Strategies UI to let the user select LIVE vs PAPER execution per strategy (and optional default per user), along with the desired poll interval.'
$ws.Range("F123").Value = 'Initial UI can be a simple dropdown or radio group for execution mode plus a select for poll intervals (e.g., 15s, 30s, 1m, 5m, 15m, 1h, 4h).'
$ws.Range("G123").Value = 'pending'
$ws.Range("I123").Value = 'Later iterations can add per-strategy overrides, presets for swing vs intraday profiles, and visual indicators of the current mode.'

# --- Row 124 ---
$ws.Range("A124").Value = 'S16'
$ws.Range("B124").Value = 'G03'
$ws.Range("C124").Value = 'Paper mode UI and analytics integration'
$ws.Range("D124").Value = 'S16_G03_TF002'
$ws.Range("E124").Value = 'Update Queue, Orders, and Analytics views to label simulated paper trades clearly and allow filtering them in/out of P&L and performance charts.'
$ws.Range("F124").Value = 'UI will likely show a small PAPER / SIM tag and default analytics to include or exclude paper trades depending on user preference.'
$ws.Range("G124").Value = 'pending'
$ws.Range("I124").Value = 'Decide default analytics behaviour (e.g., exclude paper trades by default) and expose a simple toggle in the Analytics UI.'

# --- Restore the active-cell/selection to match the edited area ---
$ws.Range("E118").Select()
